$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row/column as required by the updated TPM data.
$updates = @{
    2 = @{ G = 0.7109393333333333; H = 2.132818; I = 0.002867418779909113; J = 0.002875587087952277;
           K = 2; L = 0.6666666666666666; M = 0.638706; N = 1.916118; O = 0.5586654432763536; P = 0.5586654432763536;
           Q = 0.454081217836; R = 4.086730960524; S = 0.001601927783736865; T = 0.001606491135170618 }
    3 = @{ G = 0.7109393333333333; H = 2.132818; I = 0.002867418779909113; J = 0.002875587087952277;
           M = 0.5045649999999999; N = 1.513695; O = 0.4413345567236464; P = 0.4413345567236464;
           Q = 0.3587151047233332; R = 3.22843594251; S = 0.001265490996172247; T = 0.001269095952781659 }
    4 = @{ G = 132.4704766666667; H = 397.4114300000001; I = 0.5342907823042267; J = 0.5358127963626762;
           K = 2; L = 0.6666666666666666; M = 0.638706; N = 1.916118; O = 0.5586654432763536; P = 0.5586654432763536;
           Q = 84.60968826986002; R = 761.4871944287401; S = 0.2984897967344605; T = 0.2993400933930971 }
    5 = @{ G = 132.4704766666667; H = 397.4114300000001; I = 0.5342907823042267; J = 0.5358127963626762;
           M = 0.5045649999999999; N = 1.513695; O = 0.4413345567236464; P = 0.4413345567236464;
           Q = 66.83996605931667; R = 601.55969453385; S = 0.2358009855697661; T = 0.2364727029695791 }
    6 = @{ E = 2; F = 1; G = 2.112848; H = 4.225696; I = 0.008521711699207447; J = 0.005697324785899025;
           K = 2; L = 0.6666666666666666; M = 0.638706; N = 1.916118; O = 0.5586654432763536; P = 0.5586654432763536;
           Q = 1.349488694688; R = 8.096932168127999; S = 0.004760785843911017; T = 0.003182898477003635 }
    7 = @{ E = 2; F = 1; G = 2.112848; H = 4.225696; I = 0.008521711699207447; J = 0.005697324785899025;
           M = 0.5045649999999999; N = 1.513695; O = 0.4413345567236464; P = 0.4413345567236464;
           Q = 1.06606915112; R = 6.39641490672; S = 0.00376092585529643; T = 0.00251442630889539 }
    8 = @{ G = 112.642779; H = 337.928337; I = 0.4543200872166568; J = 0.4556142917634725;
           K = 2; L = 0.6666666666666666; M = 0.638706; N = 1.916118; O = 0.5586654432763536; P = 0.5586654432763536;
           Q = 71.945618803974; R = 647.510569235766; S = 0.2538129329142452; T = 0.2545359602710822 }
    9 = @{ G = 112.642779; H = 337.928337; I = 0.4543200872166568; J = 0.4556142917634725;
           M = 0.5045649999999999; N = 1.513695; O = 0.4413345567236464; P = 0.4413345567236464;
           Q = 56.835603786135; R = 511.520434075215; S = 0.2005071543024116; T = 0.2010783314923902 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
